$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the species-related data between rows 15 and 16, leaving the
# location/report metadata columns (C, I, P, S, T-AY) untouched.
$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr15 = "$col" + "15"
    $addr16 = "$col" + "16"
    $val15 = $ws.Range($addr15).Value2
    $val16 = $ws.Range($addr16).Value2
    $ws.Range($addr15).Value2 = $val16
    $ws.Range($addr16).Value2 = $val15
}
